$d = $word.ActiveDocument

# 1. Heading3: "Software Development and Innovation" -> "Research Leadership and Policy Impact"
$d.Content.Find.Execute(
    "Software Development and Innovation", $true, $false, $false, $false, $false,
    $true, 1, $false, "Research Leadership and Policy Impact", 2)

# 2. Bullet 1
$d.Content.Find.Execute(
    "Conceived and deployed redistricting software used by thousands of analysts nationwide",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Regular expert testimony and consultation on research methodology for journalists, elected officials, and community leaders",
    2)

# 3. Bullet 2
$d.Content.Find.Execute(
    "Developed boundary estimation system using incomplete data without ML requirements",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Research analysis used in court cases addressing housing, redistricting, and community development with rigorous methodology",
    2)

# 4. Bullet 3
$d.Content.Find.Execute(
    "Created econometric simulation platform for humanitarian intervention modeling",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Conceived and deployed cloud-based analytical software used by thousands of researchers nationwide for community-focused research",
    2)

# 5. Bullet 4
$d.Content.Find.Execute(
    "Built comprehensive survey operations platform from RFP through deployment",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Developed research frameworks and methodologies that became industry standards for community development and policy analysis",
    2)
